# Update the date line and the division problems in the worksheet table.
# Each "old" value is unique in the document, and the replacements below
# are kept in document order so that a freshly-written "new" value is
# never re-matched by a not-yet-processed "old" pattern later in this
# list (this matters for the 19÷4= / 80÷4= pair).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-09 Sunday" "2025-03-10 Monday"

Replace-Text "78÷6=" "39÷3="
Replace-Text "25÷6=" "65÷5="
Replace-Text "67÷3=" "23÷9="
Replace-Text "47÷2=" "94÷6="
Replace-Text "29÷6=" "91÷9="

Replace-Text "15÷5=" "37÷2="
Replace-Text "87÷4=" "62÷8="
Replace-Text "88÷7=" "79÷9="
Replace-Text "14÷6=" "14÷7="
Replace-Text "19÷4=" "74÷4="

Replace-Text "30÷6=" "45÷7="
Replace-Text "95÷3=" "56÷6="
Replace-Text "19÷9=" "50÷8="
Replace-Text "40÷7=" "50÷4="
Replace-Text "52÷5=" "59÷8="

Replace-Text "14÷5=" "55÷4="
Replace-Text "83÷6=" "58÷2="
Replace-Text "48÷8=" "19÷2="
Replace-Text "70÷9=" "98÷6="
Replace-Text "37÷3=" "72÷8="

Replace-Text "16÷2=" "95÷8="
Replace-Text "81÷6=" "43÷8="
Replace-Text "80÷4=" "19÷4="
Replace-Text "25÷3=" "83÷4="
Replace-Text "92÷6=" "75÷5="

Write-Output "Replacements applied"
